$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-13 Monday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-10-14 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("97-32=", $true, $true, $false, $false, $false, $true, 1, $false, "5+54=", 2) | Out-Null
$d.Content.Find.Execute("25+8=", $true, $true, $false, $false, $false, $true, 1, $false, "86-81=", 2) | Out-Null
$d.Content.Find.Execute("32+4=", $true, $true, $false, $false, $false, $true, 1, $false, "95-0=", 2) | Out-Null
$d.Content.Find.Execute("92-60=", $true, $true, $false, $false, $false, $true, 1, $false, "20-15=", 2) | Out-Null
$d.Content.Find.Execute("21+60=", $true, $true, $false, $false, $false, $true, 1, $false, "77+6=", 2) | Out-Null
$d.Content.Find.Execute("4+69=", $true, $true, $false, $false, $false, $true, 1, $false, "54-28=", 2) | Out-Null
$d.Content.Find.Execute("69-37=", $true, $true, $false, $false, $false, $true, 1, $false, "48+23=", 2) | Out-Null
$d.Content.Find.Execute("58+20=", $true, $true, $false, $false, $false, $true, 1, $false, "26+49=", 2) | Out-Null
$d.Content.Find.Execute("21-16=", $true, $true, $false, $false, $false, $true, 1, $false, "18+28=", 2) | Out-Null
$d.Content.Find.Execute("88-19=", $true, $true, $false, $false, $false, $true, 1, $false, "53+31=", 2) | Out-Null
$d.Content.Find.Execute("95-26=", $true, $true, $false, $false, $false, $true, 1, $false, "49-25=", 2) | Out-Null
$d.Content.Find.Execute("69-61=", $true, $true, $false, $false, $false, $true, 1, $false, "9+30=", 2) | Out-Null
$d.Content.Find.Execute("11+1=", $true, $true, $false, $false, $false, $true, 1, $false, "59+4=", 2) | Out-Null
$d.Content.Find.Execute("17+9=", $true, $true, $false, $false, $false, $true, 1, $false, "82-34=", 2) | Out-Null
$d.Content.Find.Execute("29-11=", $true, $true, $false, $false, $false, $true, 1, $false, "56-34=", 2) | Out-Null
$d.Content.Find.Execute("15+32=", $true, $true, $false, $false, $false, $true, 1, $false, "88-5=", 2) | Out-Null
$d.Content.Find.Execute("45-15=", $true, $true, $false, $false, $false, $true, 1, $false, "83-17=", 2) | Out-Null
$d.Content.Find.Execute("13+36=", $true, $true, $false, $false, $false, $true, 1, $false, "49+23=", 2) | Out-Null
$d.Content.Find.Execute("82-13=", $true, $true, $false, $false, $false, $true, 1, $false, "48-33=", 2) | Out-Null
$d.Content.Find.Execute("4+82=", $true, $true, $false, $false, $false, $true, 1, $false, "50+3=", 2) | Out-Null
$d.Content.Find.Execute("13+12=", $true, $true, $false, $false, $false, $true, 1, $false, "6+33=", 2) | Out-Null
$d.Content.Find.Execute("84-45=", $true, $true, $false, $false, $false, $true, 1, $false, "13+49=", 2) | Out-Null
$d.Content.Find.Execute("83-61=", $true, $true, $false, $false, $false, $true, 1, $false, "37+14=", 2) | Out-Null
$d.Content.Find.Execute("29-24=", $true, $true, $false, $false, $false, $true, 1, $false, "6+91=", 2) | Out-Null
$d.Content.Find.Execute("41+50=", $true, $true, $false, $false, $false, $true, 1, $false, "29+6=", 2) | Out-Null
$d.Content.Find.Execute("15+84=", $true, $true, $false, $false, $false, $true, 1, $false, "62+36=", 2) | Out-Null
$d.Content.Find.Execute("7+19=", $true, $true, $false, $false, $false, $true, 1, $false, "6+19=", 2) | Out-Null
$d.Content.Find.Execute("83-73=", $true, $true, $false, $false, $false, $true, 1, $false, "69-5=", 2) | Out-Null
$d.Content.Find.Execute("48-39=", $true, $true, $false, $false, $false, $true, 1, $false, "73+9=", 2) | Out-Null
$d.Content.Find.Execute("88-36=", $true, $true, $false, $false, $false, $true, 1, $false, "58-36=", 2) | Out-Null
$d.Content.Find.Execute("44-35=", $true, $true, $false, $false, $false, $true, 1, $false, "55-18=", 2) | Out-Null
$d.Content.Find.Execute("56-7=", $true, $true, $false, $false, $false, $true, 1, $false, "36-19=", 2) | Out-Null
$d.Content.Find.Execute("66-2=", $true, $true, $false, $false, $false, $true, 1, $false, "42+14=", 2) | Out-Null
$d.Content.Find.Execute("83-70=", $true, $true, $false, $false, $false, $true, 1, $false, "19-12=", 2) | Out-Null
$d.Content.Find.Execute("97-25=", $true, $true, $false, $false, $false, $true, 1, $false, "55+42=", 2) | Out-Null
$d.Content.Find.Execute("13+42=", $true, $true, $false, $false, $false, $true, 1, $false, "60+18=", 2) | Out-Null
$d.Content.Find.Execute("24-18=", $true, $true, $false, $false, $false, $true, 1, $false, "65-27=", 2) | Out-Null
$d.Content.Find.Execute("37+42=", $true, $true, $false, $false, $false, $true, 1, $false, "36+27=", 2) | Out-Null
$d.Content.Find.Execute("53+5=", $true, $true, $false, $false, $false, $true, 1, $false, "87-70=", 2) | Out-Null
$d.Content.Find.Execute("38+41=", $true, $true, $false, $false, $false, $true, 1, $false, "45-18=", 2) | Out-Null
$d.Content.Find.Execute("70+23=", $true, $true, $false, $false, $false, $true, 1, $false, "13+18=", 2) | Out-Null
$d.Content.Find.Execute("34+28=", $true, $true, $false, $false, $false, $true, 1, $false, "81-17=", 2) | Out-Null
$d.Content.Find.Execute("53+34=", $true, $true, $false, $false, $false, $true, 1, $false, "82-25=", 2) | Out-Null
$d.Content.Find.Execute("36+12=", $true, $true, $false, $false, $false, $true, 1, $false, "9+19=", 2) | Out-Null
$d.Content.Find.Execute("46-4=", $true, $true, $false, $false, $false, $true, 1, $false, "12+62=", 2) | Out-Null
$d.Content.Find.Execute("33+24=", $true, $true, $false, $false, $false, $true, 1, $false, "63-0=", 2) | Out-Null
$d.Content.Find.Execute("29+63=", $true, $true, $false, $false, $false, $true, 1, $false, "56-53=", 2) | Out-Null
$d.Content.Find.Execute("22+65=", $true, $true, $false, $false, $false, $true, 1, $false, "89-80=", 2) | Out-Null
$d.Content.Find.Execute("84-6=", $true, $true, $false, $false, $false, $true, 1, $false, "43+0=", 2) | Out-Null
$d.Content.Find.Execute("59+28=", $true, $true, $false, $false, $false, $true, 1, $false, "42-24=", 2) | Out-Null
$d.Content.Find.Execute("22-14=", $true, $true, $false, $false, $false, $true, 1, $false, "68+20=", 2) | Out-Null
$d.Content.Find.Execute("10+66=", $true, $true, $false, $false, $false, $true, 1, $false, "16+71=", 2) | Out-Null
$d.Content.Find.Execute("57+17=", $true, $true, $false, $false, $false, $true, 1, $false, "72+5=", 2) | Out-Null
$d.Content.Find.Execute("82-26=", $true, $true, $false, $false, $false, $true, 1, $false, "21+76=", 2) | Out-Null
$d.Content.Find.Execute("25+24=", $true, $true, $false, $false, $false, $true, 1, $false, "36+4=", 2) | Out-Null
$d.Content.Find.Execute("9+69=", $true, $true, $false, $false, $false, $true, 1, $false, "87-72=", 2) | Out-Null
$d.Content.Find.Execute("63-17=", $true, $true, $false, $false, $false, $true, 1, $false, "66-5=", 2) | Out-Null
$d.Content.Find.Execute("77+13=", $true, $true, $false, $false, $false, $true, 1, $false, "21-13=", 2) | Out-Null
$d.Content.Find.Execute("48+47=", $true, $true, $false, $false, $false, $true, 1, $false, "49+39=", 2) | Out-Null
$d.Content.Find.Execute("72-6=", $true, $true, $false, $false, $false, $true, 1, $false, "45+23=", 2) | Out-Null
$d.Content.Find.Execute("41+9=", $true, $true, $false, $false, $false, $true, 1, $false, "12-12=", 2) | Out-Null
$d.Content.Find.Execute("23+60=", $true, $true, $false, $false, $false, $true, 1, $false, "26+30=", 2) | Out-Null
$d.Content.Find.Execute("8+61=", $true, $true, $false, $false, $false, $true, 1, $false, "81+4=", 2) | Out-Null
$d.Content.Find.Execute("90-87=", $true, $true, $false, $false, $false, $true, 1, $false, "82-74=", 2) | Out-Null
$d.Content.Find.Execute("57+15=", $true, $true, $false, $false, $false, $true, 1, $false, "13+75=", 2) | Out-Null
$d.Content.Find.Execute("59+15=", $true, $true, $false, $false, $false, $true, 1, $false, "6+49=", 2) | Out-Null
$d.Content.Find.Execute("57+1=", $true, $true, $false, $false, $false, $true, 1, $false, "52+21=", 2) | Out-Null
$d.Content.Find.Execute("37+52=", $true, $true, $false, $false, $false, $true, 1, $false, "72-15=", 2) | Out-Null
$d.Content.Find.Execute("74-4=", $true, $true, $false, $false, $false, $true, 1, $false, "84-28=", 2) | Out-Null
$d.Content.Find.Execute("10+3=", $true, $true, $false, $false, $false, $true, 1, $false, "82-17=", 2) | Out-Null
$d.Content.Find.Execute("44+24=", $true, $true, $false, $false, $false, $true, 1, $false, "20+8=", 2) | Out-Null
$d.Content.Find.Execute("65-0=", $true, $true, $false, $false, $false, $true, 1, $false, "10+75=", 2) | Out-Null
$d.Content.Find.Execute("90-43=", $true, $true, $false, $false, $false, $true, 1, $false, "91-87=", 2) | Out-Null
$d.Content.Find.Execute("4+50=", $true, $true, $false, $false, $false, $true, 1, $false, "9+35=", 2) | Out-Null
$d.Content.Find.Execute("86-34=", $true, $true, $false, $false, $false, $true, 1, $false, "67+23=", 2) | Out-Null
$d.Content.Find.Execute("55+5=", $true, $true, $false, $false, $false, $true, 1, $false, "41+28=", 2) | Out-Null
$d.Content.Find.Execute("11+81=", $true, $true, $false, $false, $false, $true, 1, $false, "65-13=", 2) | Out-Null
$d.Content.Find.Execute("36+22=", $true, $true, $false, $false, $false, $true, 1, $false, "71-31=", 2) | Out-Null
$d.Content.Find.Execute("18+13=", $true, $true, $false, $false, $false, $true, 1, $false, "18+30=", 2) | Out-Null
$d.Content.Find.Execute("39+52=", $true, $true, $false, $false, $false, $true, 1, $false, "79-28=", 2) | Out-Null
$d.Content.Find.Execute("39+15=", $true, $true, $false, $false, $false, $true, 1, $false, "90-86=", 2) | Out-Null
$d.Content.Find.Execute("28+45=", $true, $true, $false, $false, $false, $true, 1, $false, "11-6=", 2) | Out-Null
$d.Content.Find.Execute("35+37=", $true, $true, $false, $false, $false, $true, 1, $false, "37+12=", 2) | Out-Null
$d.Content.Find.Execute("13+19=", $true, $true, $false, $false, $false, $true, 1, $false, "23-3=", 2) | Out-Null
$d.Content.Find.Execute("98-95=", $true, $true, $false, $false, $false, $true, 1, $false, "13+70=", 2) | Out-Null
$d.Content.Find.Execute("28-15=", $true, $true, $false, $false, $false, $true, 1, $false, "88-31=", 2) | Out-Null
$d.Content.Find.Execute("99-94=", $true, $true, $false, $false, $false, $true, 1, $false, "23-13=", 2) | Out-Null
$d.Content.Find.Execute("54-33=", $true, $true, $false, $false, $false, $true, 1, $false, "39-22=", 2) | Out-Null
$d.Content.Find.Execute("39+21=", $true, $true, $false, $false, $false, $true, 1, $false, "73-58=", 2) | Out-Null
$d.Content.Find.Execute("52-9=", $true, $true, $false, $false, $false, $true, 1, $false, "42+22=", 2) | Out-Null
$d.Content.Find.Execute("66-8=", $true, $true, $false, $false, $false, $true, 1, $false, "12+59=", 2) | Out-Null
$d.Content.Find.Execute("90-51=", $true, $true, $false, $false, $false, $true, 1, $false, "61-44=", 2) | Out-Null
$d.Content.Find.Execute("49-15=", $true, $true, $false, $false, $false, $true, 1, $false, "5+3=", 2) | Out-Null
$d.Content.Find.Execute("53-16=", $true, $true, $false, $false, $false, $true, 1, $false, "57+26=", 2) | Out-Null
$d.Content.Find.Execute("18+60=", $true, $true, $false, $false, $false, $true, 1, $false, "32-13=", 2) | Out-Null
$d.Content.Find.Execute("70-49=", $true, $true, $false, $false, $false, $true, 1, $false, "0+73=", 2) | Out-Null
$d.Content.Find.Execute("2+49=", $true, $true, $false, $false, $false, $true, 1, $false, "37-7=", 2) | Out-Null
$d.Content.Find.Execute("65+34=", $true, $true, $false, $false, $false, $true, 1, $false, "43+39=", 2) | Out-Null
$d.Content.Find.Execute("87-55=", $true, $true, $false, $false, $false, $true, 1, $false, "20-8=", 2) | Out-Null
$d.Content.Find.Execute("50-6=", $true, $true, $false, $false, $false, $true, 1, $false, "85-11=", 2) | Out-Null
